$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 168.85715
$ws.Range("I53").Value = 115.625
$ws.Range("J53").Value = 239.83333
$ws.Range("K53").Value = 115.625
$ws.Range("L53").Value = 239.83333
$ws.Range("M53").Value = 521.375
$ws.Range("N53").Value = -1513.83333
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents() | Out-Null
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents() | Out-Null
$ws.Range("H98").Value = 2095.4736
$ws.Range("I98").Value = 1208.4286
$ws.Range("J98").Value = 4579.2
$ws.Range("K98").Value = 1208.4286
$ws.Range("L98").Value = 4579.2
$ws.Range("M98").Value = 289.5714
$ws.Range("N98").Value = -7575.2
$ws.Range("H113").Value = 9213.138999999999
$ws.Range("I113").Value = 6746.5654
$ws.Range("J113").Value = 13577.077
$ws.Range("K113").Value = 6746.5654
$ws.Range("L113").Value = 13577.077
$ws.Range("M113").Value = -3492.5654
$ws.Range("N113").Value = -20085.077
$ws.Range("H122").Value = 2095.4736
$ws.Range("I122").Value = 1208.4286
$ws.Range("J122").Value = 4579.2
$ws.Range("K122").Value = 3625.2858
$ws.Range("L122").Value = 13737.6
$ws.Range("M122").Value = -1175.2858
$ws.Range("N122").Value = -18637.6
$ws.Range("H138").Value = 6581854
$ws.Range("I138").Value = 1623.3636
$ws.Range("J138").Value = 7695431.5
$ws.Range("K138").Value = 4870.0908
$ws.Range("L138").Value = 23086294.5
$ws.Range("M138").Value = 269.9092000000001
$ws.Range("N138").Value = -23096574.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7018.939
$ws.Range("I32").Value = 3011.561
$ws.Range("J32").Value = 27556.75
$ws.Range("K32").Value = 3011.561
$ws.Range("L32").Value = 27556.75
$ws.Range("M32").Value = -2724.561
$ws.Range("N32").Value = -28130.75
$ws.Range("H45").Value = 5449.2915
$ws.Range("I45").Value = 6191.8
$ws.Range("K45").Value = 6191.8
$ws.Range("M45").Value = -5814.8
$ws.Range("H104").Value = 74707
$ws.Range("J104").Value = 74707
$ws.Range("L104").Value = 74707
$ws.Range("N104").Value = -81695
$ws.Range("H132").Value = 3048.8394
$ws.Range("I132").Value = 2852.3618
$ws.Range("J132").Value = 4074.889
$ws.Range("K132").Value = 8557.0854
$ws.Range("L132").Value = 12224.667
$ws.Range("M132").Value = -6027.0854
$ws.Range("N132").Value = -17284.667

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 12562.5
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 12562.5
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 12562.5
$ws.Range("M50").ClearContents() | Out-Null
$ws.Range("N50").Value = -13812.5
$ws.Range("H60").Value = 18898.143
$ws.Range("I60").Value = 16146
$ws.Range("K60").Value = 16146
$ws.Range("M60").Value = -15635

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 999.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 999.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 2998.5
$ws.Range("M80").ClearContents() | Out-Null
$ws.Range("N80").Value = -4870.5
$ws.Range("H83").Value = 999.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 999.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 8995.5
$ws.Range("M83").ClearContents() | Out-Null
$ws.Range("N83").Value = -18355.5
$ws.Range("H92").Value = 1442.6
$ws.Range("J92").Value = 1442.6
$ws.Range("L92").Value = 4327.799999999999
$ws.Range("N92").Value = -6823.799999999999
$ws.Range("H132").Value = 1649.6389
$ws.Range("I132").Value = 1325.7826
$ws.Range("J132").Value = 2222.6155
$ws.Range("K132").Value = 11932.0434
$ws.Range("L132").Value = 20003.5395
$ws.Range("M132").Value = -9402.0434
$ws.Range("N132").Value = -25063.5395

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 49425
$ws.Range("J34").Value = 49425
$ws.Range("L34").Value = 49425
$ws.Range("N34").Value = -49961
$ws.Range("H70").Value = 16631.908
$ws.Range("I70").Value = 11741.667
$ws.Range("J70").Value = 22500.2
$ws.Range("K70").Value = 11741.667
$ws.Range("L70").Value = 22500.2
$ws.Range("M70").Value = -11471.667
$ws.Range("N70").Value = -23040.2
$ws.Range("H73").Value = 16631.908
$ws.Range("I73").Value = 11741.667
$ws.Range("J73").Value = 22500.2
$ws.Range("K73").Value = 11741.667
$ws.Range("L73").Value = 22500.2
$ws.Range("M73").Value = -10805.667
$ws.Range("N73").Value = -24372.2
$ws.Range("H76").Value = 49425
$ws.Range("J76").Value = 49425
$ws.Range("L76").Value = 49425
$ws.Range("N76").Value = -50055
$ws.Range("H79").Value = 49425
$ws.Range("J79").Value = 49425
$ws.Range("L79").Value = 49425
$ws.Range("N79").Value = -51609
$ws.Range("H102").Value = 45455252
$ws.Range("I102").Value = 706.35297
$ws.Range("K102").Value = 706.35297
$ws.Range("M102").Value = 915.64703
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents() | Out-Null

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3500.6553
$ws.Range("I7").Value = 3019.9614
$ws.Range("K7").Value = 3019.9614
$ws.Range("M7").Value = -2907.9614
$ws.Range("H9").Value = 379.5
$ws.Range("I9").Value = 202.33333
$ws.Range("K9").Value = 202.33333
$ws.Range("M9").Value = 21.66667000000001
$ws.Range("H11").Value = 6007
$ws.Range("J11").Value = 6007
$ws.Range("L11").Value = 6007
$ws.Range("N11").Value = -6287
$ws.Range("H13").Value = 5327.0586
$ws.Range("I13").Value = 5327.0586
$ws.Range("K13").Value = 5327.0586
$ws.Range("M13").Value = -5187.0586
$ws.Range("H22").Value = 2018.25
$ws.Range("J22").Value = 2441.2144
$ws.Range("L22").Value = 2441.2144
$ws.Range("N22").Value = -3031.2144
$ws.Range("H25").Value = 1426.6342
$ws.Range("I25").Value = 1492.8572
$ws.Range("J25").Value = 1040.3334
$ws.Range("K25").Value = 1492.8572
$ws.Range("L25").Value = 1040.3334
$ws.Range("M25").Value = -1262.8572
$ws.Range("N25").Value = -1500.3334
$ws.Range("H27").Value = 2018.25
$ws.Range("J27").Value = 2441.2144
$ws.Range("L27").Value = 2441.2144
$ws.Range("N27").Value = -2655.2144
$ws.Range("H46").Value = 1045.3
$ws.Range("J46").Value = 1362.5
$ws.Range("L46").Value = 1362.5
$ws.Range("N46").Value = -1738.5
$ws.Range("H81").Value = 29181
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 29181
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 29181
$ws.Range("M81").ClearContents() | Out-Null
$ws.Range("N81").Value = -31177
$ws.Range("H82").Value = 9525.538
$ws.Range("J82").Value = 2400
$ws.Range("L82").Value = 2400
$ws.Range("N82").Value = -3122
$ws.Range("H84").Value = 29181
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 29181
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 87543
$ws.Range("M84").ClearContents() | Out-Null
$ws.Range("N84").Value = -97527
$ws.Range("H85").Value = 9525.538
$ws.Range("J85").Value = 2400
$ws.Range("L85").Value = 2400
$ws.Range("N85").Value = -4896
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents() | Out-Null
$ws.Range("H100").Value = 2029.6
$ws.Range("J100").Value = 2698
$ws.Range("L100").Value = 2698
$ws.Range("N100").Value = -3780
$ws.Range("H126").Value = 3500.6553
$ws.Range("I126").Value = 3019.9614
$ws.Range("K126").Value = 9059.8842
$ws.Range("M126").Value = -6589.8842

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 22222
$ws.Range("J76").Value = 22222
$ws.Range("L76").Value = 22222
$ws.Range("N76").Value = -22852
$ws.Range("H79").Value = 22222
$ws.Range("J79").Value = 22222
$ws.Range("L79").Value = 22222
$ws.Range("N79").Value = -24406
$ws.Range("H105").Value = 13958.125
$ws.Range("J105").Value = 13958.125
$ws.Range("L105").Value = 13958.125
$ws.Range("N105").Value = -20946.125
$ws.Range("H126").Value = 14988.375
$ws.Range("I126").Value = 16915.285
$ws.Range("K126").Value = 50745.855
$ws.Range("M126").Value = -48275.855
